$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> FAPs, Ngf -> Sorcs3)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.814087666666667
$ws.Range("H2").Value = 5.442263
$ws.Range("I2").Value = 0.1211063206477811
$ws.Range("J2").Value = 0.1211063206477811
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005316
$ws.Range("N2").Value = 0.015948
$ws.Range("Q2").Value = 0.009643690036
$ws.Range("R2").Value = 0.086793210324
$ws.Range("S2").Value = 0.1211063206477811
$ws.Range("T2").Value = 0.1211063206477811

# Row 3 (FAPs)
$ws.Range("I3").Value = 0.3837539427192561
$ws.Range("J3").Value = 0.3837539427192561
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.005316
$ws.Range("N3").Value = 0.015948
$ws.Range("Q3").Value = 0.03055830656800001
$ws.Range("R3").Value = 0.275024759112
$ws.Range("S3").Value = 0.3837539427192561
$ws.Range("T3").Value = 0.3837539427192561

# Row 4 (MuSCs)
$ws.Range("G4").Value = 7.416845666666667
$ws.Range("H4").Value = 22.250537
$ws.Range("I4").Value = 0.4951397366329628
$ws.Range("J4").Value = 0.4951397366329628
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.005316
$ws.Range("N4").Value = 0.015948
$ws.Range("Q4").Value = 0.039427951564
$ws.Range("R4").Value = 0.354851564076
$ws.Range("S4").Value = 0.4951397366329628
$ws.Range("T4").Value = 0.4951397366329628
